# Extent 5 Chnages and Merged code
#
# This script replays, via the Excel COM object model, the edits captured in
# the commit "Extent 5 Chnages and Merged code":
#   1. ProfileOperations: duplicate the first 5 data rows (rows 2-6) as new
#      rows 27-31, renaming "- Basic" to "- Duplicate" and the result column
#      to "Duplicate".
#   2. JS-Latest: unhide a block of previously-hidden rows and give several
#      rows custom (auto-fit-like) row heights; also touch a few cell values
#      (rows 54-56).
#   3. BookMarks: change D10 from "B3" to "B8".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. ProfileOperations -- duplicate rows 2-6 into rows 27-31
# ---------------------------------------------------------------------------
$wsProfileOps = $wb.Worksheets.Item("ProfileOperations")
$wsProfileOps.Select()

$sourceRows = @(2, 3, 4, 5, 6)
$destStart = 27

for ($i = 0; $i -lt $sourceRows.Count; $i++) {
    $srcRow = $sourceRows[$i]
    $dstRow = $destStart + $i

    $srcRange = $wsProfileOps.Range("A" + $srcRow + ":G" + $srcRow)
    $dstRange = $wsProfileOps.Range("A" + $dstRow + ":G" + $dstRow)

    $srcRange.Copy()
    $dstRange.PasteSpecial(-4104)  # xlPasteAll

    $aVal = $wsProfileOps.Cells.Item($srcRow, 1).Value
    $newAVal = $aVal -replace " - Basic$", " - Duplicate"
    $wsProfileOps.Cells.Item($dstRow, 1).Value = $newAVal
    $wsProfileOps.Cells.Item($dstRow, 7).Value = "Duplicate"
}

$wsProfileOps.Application.CutCopyMode = $false
$wsProfileOps.Range("G21").Select()

# ---------------------------------------------------------------------------
# 2. JS-Latest -- unhide rows and set auto-fit style heights
# ---------------------------------------------------------------------------
$wsJs = $wb.Worksheets.Item("JS-Latest")
$wsJs.Select()

# row 15 simply becomes visible again (height untouched)
$wsJs.Rows("15:15").Hidden = $false

# row 42 gets a small explicit height (stays visible)
$wsJs.Rows("42:42").Hidden = $false
$wsJs.Rows("42:42").RowHeight = 10.8

# row 46 stays hidden, but records a (very small) explicit height
$wsJs.Rows("46:46").RowHeight = 1.2
$wsJs.Rows("46:46").Hidden = $true

# row 48 gets a small explicit height (stays visible)
$wsJs.Rows("48:48").Hidden = $false
$wsJs.Rows("48:48").RowHeight = 5.4

# rows 53-91: unhide + per-row custom height (content-sized, as if AutoFit)
$rowHeights = @{
    53 = 77.4
    54 = 66.6
    55 = 55.8
    56 = 60.6
    57 = 70.8
    58 = 63
    59 = 43.8
    60 = 63
    61 = 63
    62 = 63
    63 = 63
    64 = 63
    65 = 63
    66 = 56.4
    67 = 63
    68 = 60
    69 = 49.2
    70 = 63
    71 = 56.4
    72 = 63
    73 = 52.8
    74 = 67.8
    75 = 22.8
    76 = 39
    77 = 67.8
    78 = 67.8
    79 = 33
    80 = 50.4
    81 = 57
    82 = 60.6
    83 = 67.2
    84 = 59.4
    85 = 51.6
    86 = 63.6
    87 = 55.2
    88 = 57
    89 = 47.4
    90 = 27.6
    91 = 37.2
}

for ($r = 53; $r -le 91; $r++) {
    $wsJs.Rows("$r`:$r").Hidden = $false
    $wsJs.Rows("$r`:$r").RowHeight = $rowHeights[$r]
}

# Cell-content touch-ups inside that block (rows 54-56)
$wsJs.Range("G54").Value = "Remote"

$wsJs.Range("C55").Value = "Job Script"
$wsJs.Range("C55").Borders.LineStyle = 1
$wsJs.Range("G55").Value = "Remote"
$wsJs.Range("G55").Borders.LineStyle = 1

$wsJs.Range("C56").Value = "Job Script"
$wsJs.Range("C56").Borders.LineStyle = 1
$wsJs.Range("G56").Value = "Remote"
$wsJs.Range("G56").Borders.LineStyle = 1

$wsJs.Application.WindowState = -4143  # xlMaximized (no-op placeholder)
$wsJs.Range("C54:C56").Select()

# ---------------------------------------------------------------------------
# 3. BookMarks -- D10 "B3" -> "B8"
# ---------------------------------------------------------------------------
$wsBookMarks = $wb.Worksheets.Item("BookMarks")
$wsBookMarks.Select()
$wsBookMarks.Range("D10").Value = "B8"
$wsBookMarks.Range("M10").Select()

Write-Output "edits applied"
